$d = $word.ActiveDocument

# The document currently ends with an empty paragraph (the 7th paragraph).
# We replace its content and append a new paragraph after it, matching the
# author's newly typed text: "Testtesttest " then a new paragraph "Msd".

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range

$r.Font.Name = "Segoe UI"
$r.Font.Color = 1115921        # BGR for 101214 -> 0x141210 = 1115921
$r.Font.Size = 10.5
$r.Shading.BackgroundPatternColor = 16645372   # BGR for FCFDFE -> 0xFEFDFC

$r.InsertAfter("Testtesttest ")
$r.InsertParagraphAfter()

$d.Paragraphs.Last.Range.Text = "Msd"
$d.Paragraphs.Last.Range.Font.Name = "Segoe UI"
$d.Paragraphs.Last.Range.Font.Color = 1115921
$d.Paragraphs.Last.Range.Font.Size = 10.5
$d.Paragraphs.Last.Range.Shading.BackgroundPatternColor = 16645372
